$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (row 6 and row 7) to the feed logs sheet
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "2024-06-15 00:57:12"
$ws.Range("D6").Value = 200
$ws.Range("E6").Value = 9

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "2024-06-15 00:57:13"
$ws.Range("D7").Value = 200
$ws.Range("E7").Value = 0
